$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "dd"
$ws.Range("C2").Value = "dd"
$ws.Range("D2").Value = "dd"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "Driver"
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = $true

$ws.Columns.AutoFit()
